$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Sheet, $Addr, $Val)
    $r = $Sheet.Range($Addr)
    $origStyle = $r.Style
    $r.NumberFormat = "@"
    $r.Value = $Val
    $r.Style = $origStyle
}

Set-TextValue $ws 'D2' '67.197.33'
Set-TextValue $ws 'D3' '3.497.13'
Set-TextValue $ws 'E3' '  -4.27%  '
Set-TextValue $ws 'E4' '  -0.13%  '
Set-TextValue $ws 'D5' '200.59'
Set-TextValue $ws 'E5' '  +2.15%  '
Set-TextValue $ws 'D6' '550.44'
Set-TextValue $ws 'E6' '  -5.33%  '
Set-TextValue $ws 'D7' '3.488.29'
Set-TextValue $ws 'E7' '  -4.41%  '
Set-TextValue $ws 'D8' '0.605'
Set-TextValue $ws 'E8' '  -2.66%  '
Set-TextValue $ws 'E9' '  +0.04%  '
Set-TextValue $ws 'D10' '0.652'
Set-TextValue $ws 'E10' '  -4.69%  '
Set-TextValue $ws 'D11' '62.91'
Set-TextValue $ws 'E11' '  +9.11%  '
Set-TextValue $ws 'D12' '0.143'
Set-TextValue $ws 'E12' '  -7.57%  '
Set-TextValue $ws 'D13' '0.0000270'
Set-TextValue $ws 'E13' '  -9.26%  '
Set-TextValue $ws 'D14' '9.79'
Set-TextValue $ws 'E14' '  -4.27%  '
Set-TextValue $ws 'D15' '4.046.03'
Set-TextValue $ws 'E15' '  -4.50%  '
Set-TextValue $ws 'D16' '3.481.60'
Set-TextValue $ws 'E16' '  -4.75%  '
Set-TextValue $ws 'D17' '0.124'
Set-TextValue $ws 'E17' '  -2.14%  '
Set-TextValue $ws 'B18' 'Chainlink'
Set-TextValue $ws 'C18' 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-TextValue $ws 'D18' '18.30'
Set-TextValue $ws 'E18' '  -2.23%  '
Set-TextValue $ws 'B19' 'WrappedBTC'
Set-TextValue $ws 'C19' 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
Set-TextValue $ws 'D19' '66.867.66'
Set-TextValue $ws 'D20' '11.77'
Set-TextValue $ws 'E20' '  -6.89%  '
Set-TextValue $ws 'E21' '  -6.33%  '
Set-TextValue $ws 'D22' '390.50'
Set-TextValue $ws 'E22' '  -3.33%  '
Set-TextValue $ws 'D23' '12.47'
Set-TextValue $ws 'E23' '  -3.32%  '
Set-TextValue $ws 'D24' '3.98'
Set-TextValue $ws 'E24' '  -6.82%  '
Set-TextValue $ws 'D25' '82.50'
Set-TextValue $ws 'E25' '  -4.53%  '
Set-TextValue $ws 'D26' '3.87'
Set-TextValue $ws 'E26' '  -0.18%  '
Set-TextValue $ws 'D27' '2.81'
Set-TextValue $ws 'E27' '  -5.46%  '
Set-TextValue $ws 'D28' '12.17'
Set-TextValue $ws 'E28' '  -4.29%  '
Set-TextValue $ws 'D29' '8.80'
Set-TextValue $ws 'E29' '  -4.75%  '
Set-TextValue $ws 'D30' '30.91'
Set-TextValue $ws 'E30' '  -3.05%  '
Set-TextValue $ws 'D31' '675.76'
Set-TextValue $ws 'E31' '  -2.22%  '
Set-TextValue $ws 'D32' '6.92'
Set-TextValue $ws 'E32' '  -15.79%  '
Set-TextValue $ws 'D33' '11.67'
Set-TextValue $ws 'E33' '  -5.07%  '
Set-TextValue $ws 'D34' '63.65'
Set-TextValue $ws 'E34' '  -2.14%  '
Set-TextValue $ws 'E35' '  -7.04%  '
Set-TextValue $ws 'D36' '38.69'
Set-TextValue $ws 'E36' '  -9.80%  '
Set-TextValue $ws 'E37' '  +0.04%  '
Set-TextValue $ws 'D38' '0.395'
Set-TextValue $ws 'E38' '  -7.53%  '
Set-TextValue $ws 'B39' 'Kaspa'
Set-TextValue $ws 'C39' 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextValue $ws 'D39' '0.131'
Set-TextValue $ws 'E39' '  -4.56%  '
Set-TextValue $ws 'B40' 'FirstDigitalUSD'
Set-TextValue $ws 'C40' 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
Set-TextValue $ws 'D40' '0.997'
Set-TextValue $ws 'E40' '  -0.22%  '
Set-TextValue $ws 'D41' '3.058.42'
Set-TextValue $ws 'E41' '  -5.63%  '
Set-TextValue $ws 'D42' '2.96'
Set-TextValue $ws 'E42' '  -5.66%  '
Set-TextValue $ws 'B43' 'Fetch.AI'
Set-TextValue $ws 'C43' 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
Set-TextValue $ws 'D43' '2.59'
Set-TextValue $ws 'E43' '  -9.47%  '
Set-TextValue $ws 'B44' 'PEPE'
Set-TextValue $ws 'C44' 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
Set-TextValue $ws 'D44' '0.0₃0669'
Set-TextValue $ws 'E44' '  -16.09%  '
Set-TextValue $ws 'D45' '2.77'
Set-TextValue $ws 'E45' '  +5.31%  '
Set-TextValue $ws 'B46' 'dogwifhat'
Set-TextValue $ws 'C46' 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
Set-TextValue $ws 'D46' '2.72'
Set-TextValue $ws 'E46' '  -7.93%  '
Set-TextValue $ws 'B47' 'VeChain'
Set-TextValue $ws 'C47' 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue $ws 'D47' '0.0397'
Set-TextValue $ws 'E47' '  -6.63%  '
Set-TextValue $ws 'D48' '0.127'
Set-TextValue $ws 'E48' '  -4.61%  '
Set-TextValue $ws 'D49' '137.02'
Set-TextValue $ws 'E49' '  -4.60%  '
Set-TextValue $ws 'D50' '8.19'
Set-TextValue $ws 'E50' '  -8.33%  '
Set-TextValue $ws 'D51' '2.85'
Set-TextValue $ws 'E51' '  -8.85%  '
